# "color palette: add new based on existing" task is marked Done:
# it moves from the "Active" sheet (Todo list) to the top of the
# "Inactive" sheet (Done list), with a Done date recorded.

$wb = $excel.ActiveWorkbook

$active = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")

# Capture the task's data before removing it from Active (row 5).
$taskId = $active.Range("A5").Value2
$taskTitle = $active.Range("B5").Value2
$taskCategory = $active.Range("D5").Value2
$taskCreated = $active.Range("E5").Value2

# Remove the row from Active; rows below shift up.
$active.Rows.Item(5).Delete()

# Insert a new row at the top of the data in Inactive; rows below shift down.
$inactive.Rows.Item(2).Insert()

$inactive.Range("A2").Value = $taskId
$inactive.Range("B2").Value = $taskTitle
$inactive.Range("C2").Value = "Done"
$inactive.Range("D2").Value = $taskCategory
# Created/Done columns hold plain "m/d/yyyy" text (shared strings), not real
# dates - prefix with an apostrophe so Excel stores literal text instead of
# auto-converting to a date serial.
$inactive.Range("E2").Value = "'" + $taskCreated
$inactive.Range("F2").Value = "'" + $taskCreated

# The inserted row picked up the bold header formatting (and the forced text
# format above); reset the whole row back to the plain default style used by
# the rest of the data rows.
$inactive.Range("A2:F2").Style = "Normal"
